$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto symbol price/volume data (GitHub Actions scrape refresh)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "256.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.38%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.217"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.30%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05857"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.90%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.698"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.76%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8705"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.68%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9572"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "11.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1411"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.99%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07161"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03210"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09210"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.49%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001552"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.07%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006091"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.78%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005921"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.53%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.212"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.47%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.49%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3179"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.69%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03460"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.30%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.533"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.46%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.68%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001229"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.96%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004564"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.31%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001201"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.19%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001466"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "1.29%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03816"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005622"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.76%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1102"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.87%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.41%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009823"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.73%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005392"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.31%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.14%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.09001"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "11.37%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002130"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-24.89%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.14%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"
